$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New transaction rows appended by this sync batch (rows 328-374).
# Column B holds numeric-looking account/phone numbers that must stay
# text (matches the source data, which stores them as inline strings),
# so the destination range is pre-formatted as Text before the values
# are written - otherwise Excel auto-converts the digit strings to numbers.
$ws.Range("B328:B374").NumberFormat = "@"

$rows = @(
    @('2026-02-15 12:32:43', '237674895877', 'NZOUNDI NGUEDJO EDITH TOP MOBIL TELECOM', 1326911),
    @('2026-02-15 15:07:02', '237674956331', 'PAULIN POUALEU', 168225),
    @('2026-02-15 13:18:47', '237674979451', 'EMPIRE COMPANY   LIMITED SWKBA CBOX R4 TCHOUALA GLADIS NADEGE', 377429),
    @('2026-02-15 12:29:08', '237675329176', 'NGHOMBOMBONG ALPHONSE NKWANUI LA NEGRESSE SARL', 204877),
    @('2026-02-15 01:35:22', '237675889982', 'LAIC YDRAN LIEDJI DAYA', 182245),
    @('2026-02-15 10:52:48', '237676076193', 'AMADOU ROBERT ETS MOBILE FINANCIAL SERVICES MFS', 602610),
    @('2026-02-12 14:11:34', '237677567788', 'CARINE DEVILLIERS ESSOH', 448951),
    @('2026-02-14 10:49:21', '237678722616', 'KWEDI BELLE GEORGES GAEL ETS MOBILE FINANCIAL SERVICES MFS', 6039),
    @('2026-02-15 16:03:35', '237679580678', 'Ngo Ndjeng Marie Pauline', 5754),
    @('2026-02-15 11:41:35', '237680116452', 'LA NEGRESSE SARL NDE HENRI', 97795),
    @('2026-02-15 10:06:37', '237680415220', 'TIPA CHRISTELLE NADIA LA NEGRESSE SARL', 10262),
    @('2026-02-14 13:55:46', '237681603496', 'MAGNI YANSI DORIS ETS MOBILE FINANCIAL SERVICES MFS', 130572),
    @('2026-02-15 14:57:21', '237681657461', 'LA NEGRESSE SARL MOKAM KOUAM VIVIANE', 70335),
    @('2026-02-15 13:45:49', '237681662596', 'LA NEGRESSE SARL LTDLA CBOX R1 TIOKENG SANDRINE', 194688),
    @('2026-02-15 09:38:38', '237681677617', 'ETS MOBILE FINANCIAL SERVICES MFS MANFOUO TCHOUALA HUGUETTE', 27208),
    @('2026-02-15 13:44:32', '237681679096', 'TAMNOU NGANGO ULRICH BERNARD ETS MOBILE FINANCIAL SERVICES MFS', 97384),
    @('2026-02-15 16:43:25', '237683279255', 'Deuffi Yvonne laurette LENA GLOBAL', 34169),
    @('2026-02-14 16:04:50', '237683353137', 'VOUGMO NGUEMO MERLIN WILLIAM ETS MOBILE FINANCIAL SERVICES MFS', 318),
    @('2026-02-15 11:35:46', '237683396173', 'ENOMA NDJAH PAULINE SONIA ETS MOBILE FINANCIAL SERVICES MFS', 223792),
    @('2026-02-15 13:49:14', '237683557193', 'KOUYEKE MONIQUE LA NEGRESSE SARL', 254),
    @('2026-02-15 16:23:31', '237671646117', 'FLORENCE MEDONGUE', 516490),
    @('2026-02-15 12:10:52', '237674440808', 'HOUMI EPSE MVEINGUE GUINDOP TATIANA ETS LE CONTENT', 22158),
    @('2026-02-15 16:37:06', '237675788721', 'ETS MOBILE FINANCIAL SERVICES MFS TCHAWE MBOUGA JUDITH FLORE', 513322),
    @('2026-02-15 12:22:44', '237676036914', 'NDENGUE ELOKO SAMUEL HERVE ETS MOBILE FINANCIAL SERVICES MFS', 201711),
    @('2026-02-15 12:08:24', '237679550294', 'NGOUNOU JOSIANE GAELLE CHIC MOBILE SARL', 8438),
    @('2026-02-15 09:21:13', '237679604574', 'KINGUE KOMBI VICTORINE SIDONIE VISION TRADING COMPLEX AND TECHNOLOGIES SARL VISION TRADING COMPLEX', 16343),
    @('2026-02-15 11:43:46', '237682370358', 'TATALONG NDAH DARYL', 270541),
    @('2026-02-15 11:40:17', '237682639044', 'NJIMEGNIE LINDA NANA', 568123),
    @('2026-02-15 13:56:24', '237683023087', 'FAGHUIE ABIBA', 342523),
    @('2026-02-15 14:31:34', '237650874464', 'ASSONFACK VANESSA ATB POINT COM', 471335),
    @('2026-02-15 09:41:32', '237653854849', 'Jean Boreau Holl', 51705),
    @('2026-02-15 16:45:52', '237654164073', 'NGUIDJOL SIMONE ASTRIDE SPECTRUM SPECTRUM', 52578),
    @('2026-02-15 16:31:02', '237671357520', 'CALICE WOTI EPSE DJOMO', 22850),
    @('2026-02-13 13:41:17', '237672587687', 'Kamaha Tomy Nadine LA NEGRESSE SARL', 247714),
    @('2026-02-15 14:40:50', '237674240552', 'AUGUSTINE NGO BAYOI', 60758),
    @('2026-02-15 14:11:21', '237675239360', 'NGO NKOT MARIE LOUISE ETS MOBILE FINANCIAL SERVICES MFS', 118239),
    @('2026-02-15 15:17:14', '237675626141', 'FLORENCE NGUEFACK', 52169),
    @('2026-02-15 19:24:24', '237676840777', 'ETP109 ETP', 994075),
    @('2026-02-15 17:11:12', '237677833877', 'MEDJOM TAGNE MICHELLE GUILENE LA NEGRESSE SARL', 144763),
    @('2026-02-15 10:02:37', '237678854978', 'rostel kenzoh gignole', 250763),
    @('2026-02-15 14:55:22', '237679422591', 'ETS LE CONTENT 42', 403900),
    @('2026-02-15 14:13:01', '237650353920', 'MENIAPI HELENE EDOSSINE TOP MOBIL TELECOM', 1305524),
    @('2026-02-15 10:59:08', '237651927448', 'charity aben awalah', 28091),
    @('2026-02-14 12:20:24', '237653294562', 'NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS', 271672),
    @('2026-02-14 13:53:49', '237678046498', 'MFS SOCAVER', 303),
    @('2026-02-15 14:55:21', '237679428698', 'ETS LE CONTENT 29', 13),
    @('2026-02-14 17:33:53', '237679551262', 'LA NEGRESSE LTDLA CBOX R1 MEGNE JUDITH', 17349)
)

$startRow = 328
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
}
